$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 614.7826
$ws.Range("I28").Value = 588.5454999999999
$ws.Range("J28").Value = 638.8333
$ws.Range("K28").Value = 588.5454999999999
$ws.Range("L28").Value = 638.8333
$ws.Range("M28").Value = -103.5454999999999
$ws.Range("N28").Value = -1608.8333

$ws.Range("H98").Value = 26558.742
$ws.Range("I98").Value = 952
$ws.Range("J98").Value = 89153
$ws.Range("K98").Value = 952
$ws.Range("L98").Value = 89153
$ws.Range("M98").Value = 546
$ws.Range("N98").Value = -92149

$ws.Range("H107").Value = 605.8946999999999
$ws.Range("I107").Value = 506.22223
$ws.Range("K107").Value = 506.22223
$ws.Range("M107").Value = 1413.77777

$ws.Range("H121").Value = 2573.3333
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 2573.3333
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 7719.999899999999
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -11213.9999

$ws.Range("H122").Value = 26558.742
$ws.Range("I122").Value = 952
$ws.Range("J122").Value = 89153
$ws.Range("K122").Value = 2856
$ws.Range("L122").Value = 267459
$ws.Range("M122").Value = -406
$ws.Range("N122").Value = -272359

$ws.Range("H129").Value = 1088.841
$ws.Range("I129").Value = 1091.8
$ws.Range("J129").Value = 1087.3103
$ws.Range("K129").Value = 3275.4
$ws.Range("L129").Value = 3261.9309
$ws.Range("M129").Value = 1724.6
$ws.Range("N129").Value = -13261.9309

$ws.Range("H137").Value = 3210635.5
$ws.Range("I137").Value = 6994032
$ws.Range("K137").Value = 20982096
$ws.Range("M137").Value = -20979546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1453.6046
$ws.Range("I2").Value = 1524.75
$ws.Range("J2").Value = 1246.6364
$ws.Range("K2").Value = 1524.75
$ws.Range("L2").Value = 1246.6364
$ws.Range("M2").Value = -1411.75
$ws.Range("N2").Value = -1472.6364

$ws.Range("H61").Value = 2697.853
$ws.Range("I61").Value = 1790
$ws.Range("J61").Value = 3259.8572
$ws.Range("K61").Value = 1790
$ws.Range("L61").Value = 3259.8572
$ws.Range("M61").Value = -1578
$ws.Range("N61").Value = -3683.8572

$ws.Range("H116").Value = 1453.6046
$ws.Range("I116").Value = 1524.75
$ws.Range("J116").Value = 1246.6364
$ws.Range("K116").Value = 1524.75
$ws.Range("L116").Value = 1246.6364
$ws.Range("M116").Value = 769.25
$ws.Range("N116").Value = -5834.6364

$ws.Range("H122").Value = 1741.4615
$ws.Range("I122").Value = 1589.4762
$ws.Range("J122").Value = 2379.8
$ws.Range("K122").Value = 4768.4286
$ws.Range("L122").Value = 7139.400000000001
$ws.Range("M122").Value = -2318.4286
$ws.Range("N122").Value = -12039.4

$ws.Range("H136").Value = 2697.853
$ws.Range("I136").Value = 1790
$ws.Range("J136").Value = 3259.8572
$ws.Range("K136").Value = 5370
$ws.Range("L136").Value = 9779.571599999999
$ws.Range("M136").Value = -2820
$ws.Range("N136").Value = -14879.5716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1453.6046
$ws.Range("I3").Value = 1524.75
$ws.Range("J3").Value = 1246.6364
$ws.Range("K3").Value = 1524.75
$ws.Range("L3").Value = 1246.6364
$ws.Range("M3").Value = -1410.75
$ws.Range("N3").Value = -1474.6364

$ws.Range("H105").Value = 4149.3335
$ws.Range("I105").Value = 3702
$ws.Range("J105").Value = 4468.857
$ws.Range("K105").Value = 3702
$ws.Range("L105").Value = 4468.857
$ws.Range("M105").Value = -1955
$ws.Range("N105").Value = -7962.857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1647
$ws.Range("I16").Value = 1647
$ws.Range("K16").Value = 1647
$ws.Range("M16").Value = -1360

$ws.Range("H106").Value = 52571.43
$ws.Range("J106").Value = 52571.43
$ws.Range("L106").Value = 52571.43
$ws.Range("N106").Value = -55095.43

$ws.Range("H107").Value = 685.5484
$ws.Range("I107").Value = 617.875
$ws.Range("J107").Value = 757.73334
$ws.Range("K107").Value = 617.875
$ws.Range("L107").Value = 757.73334
$ws.Range("M107").Value = 1302.125
$ws.Range("N107").Value = -4597.73334

$ws.Range("H113").Value = 1647
$ws.Range("I113").Value = 1647
$ws.Range("K113").Value = 1647
$ws.Range("M113").Value = 523

$ws.Range("H132").Value = 1080234.5
$ws.Range("I132").Value = 2603.2856
$ws.Range("J132").Value = 2337470.8
$ws.Range("K132").Value = 7809.8568
$ws.Range("L132").Value = 7012412.399999999
$ws.Range("M132").Value = -5279.8568
$ws.Range("N132").Value = -7017472.399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 430.18182
$ws.Range("I23").Value = 392
$ws.Range("J23").Value = 452
$ws.Range("K23").Value = 1176
$ws.Range("L23").Value = 1356
$ws.Range("M23").Value = -941
$ws.Range("N23").Value = -1826

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 3215.8462
$ws.Range("I107").Value = 589
$ws.Range("K107").Value = 589
$ws.Range("M107").Value = 1331

$ws.Range("H113").Value = 1714.2858
$ws.Range("I113").Value = 2060
$ws.Range("J113").Value = 850
$ws.Range("K113").Value = 2060
$ws.Range("L113").Value = 850
$ws.Range("M113").Value = 110
$ws.Range("N113").Value = -5190

$ws.Range("H126").Value = 20223
$ws.Range("I126").Value = 37008
$ws.Range("J126").Value = 3438
$ws.Range("K126").Value = 111024
$ws.Range("L126").Value = 10314
$ws.Range("M126").Value = -108554
$ws.Range("N126").Value = -15254

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1808.1936
$ws.Range("I136").Value = 1454.4
$ws.Range("J136").Value = 3282.3333
$ws.Range("K136").Value = 4363.200000000001
$ws.Range("L136").Value = 9846.999899999999
$ws.Range("M136").Value = -1813.200000000001
$ws.Range("N136").Value = -14946.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 6251406
$ws.Range("I107").Value = 1150
$ws.Range("J107").Value = 16668500
$ws.Range("K107").Value = 3450
$ws.Range("L107").Value = 50005500
$ws.Range("M107").Value = -1530
$ws.Range("N107").Value = -50009340

$ws.Range("H132").Value = 3108913.2
$ws.Range("I132").Value = 6213635.5
$ws.Range("J132").Value = 4190.857
$ws.Range("K132").Value = 18640906.5
$ws.Range("L132").Value = 12572.571
$ws.Range("M132").Value = -18638376.5
$ws.Range("N132").Value = -17632.571

$ws.Range("H136").Value = 460869.47
$ws.Range("I136").Value = 534236.75
$ws.Range("J136").Value = 2324.25
$ws.Range("K136").Value = 1602710.25
$ws.Range("L136").Value = 6972.75
$ws.Range("M136").Value = -1600160.25
$ws.Range("N136").Value = -12072.75
